$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.01514828764759746, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1, 0.9822431866464301)
    3 = @(0.003994804209775715, 1.667794583268128, 3.900430680208489, 8.660232485948974, 1, 14.23245255363537)
    4 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    5 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 0, 9.906394106792828)
    6 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    7 = @(0.01514828764759746, 0.00007097389502863649, 0.1575252929769615, 0.496779210170732, 0, 0.6695237646903196)
    8 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 4.429675500412797)
    9 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
